$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.270.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.497.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.525"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.72%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("E14").Value = "  +1.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.887.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.498.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.201.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.19%  "
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("E20").Value = "  +3.26%  "
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("E22").Value = "  +13.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "247.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.52%  "
$ws.Range("E31").Value = "  +8.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("E35").Value = "  +2.95%  "
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("E37").Value = "  +4.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.60%  "
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("E40").Value = "  +1.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "121.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.98%  "
$ws.Range("E44").Value = "  +2.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.995.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("E46").Value = "  +2.74%  "
$ws.Range("E47").Value = "  -1.70%  "
$ws.Range("E48").Value = "  -3.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.05%  "
